# Update NATMI LR-pair TPM-derived metrics (Ntng2-Lrrc4) with new TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> ordered list of (column, new value) pairs to write.
$updates = @{
    2 = @{ "G" = 8.031680333333334; "H" = 24.095041; "I" = 0.353073148880341; "J" = 0.3600143049225579; "M" = 1.546825; "N" = 4.640475; "O" = 0.09335776776421496; "P" = 0.09839297464651384; "Q" = 12.42360393160834; "R" = 111.812435384475; "S" = 0.03296212103695097; "T" = 0.03542287837662754 }
    3 = @{ "G" = 8.031680333333334; "H" = 24.095041; "I" = 0.353073148880341; "J" = 0.3600143049225579; "O" = 0.3226929381796237; "P" = 0.3400972286002662; "Q" = 42.94242837506356; "R" = 386.4818553755721; "S" = 0.113934211804529; "T" = 0.1224398673606131 }
    4 = @{ "G" = 8.031680333333334; "H" = 24.095041; "I" = 0.353073148880341; "J" = 0.3600143049225579; "M" = 3.750299666666667; "N" = 11.250899; "O" = 0.2263472631531553; "P" = 0.238555195331833; "Q" = 30.12120807687323; "R" = 271.090872691859; "S" = 0.07991714094193171; "T" = 0.08588328283305489 }
    5 = @{ "G" = 8.031680333333334; "H" = 24.095041; "I" = 0.353073148880341; "J" = 0.3600143049225579; "M" = 2.543696; "N" = 5.087391999999999; "O" = 0.1535233658822184; "P" = 0.1078690504900635; "Q" = 20.43015313717867; "R" = 122.580918823072; "S" = 0.05420497821874356; "T" = 0.0388344012348365 }
    6 = @{ "G" = 8.031680333333334; "H" = 24.095041; "I" = 0.353073148880341; "J" = 0.3600143049225579; "M" = 3.381336; "N" = 10.144008; "O" = 0.2040786650207874; "P" = 0.2150855509313235; "Q" = 27.157809851592; "R" = 244.420288664328; "S" = 0.07205469687818571; "T" = 0.07743387511742585 }
    7 = @{ "I" = 0.1399128257358017; "J" = 0.1426634080409682; "M" = 1.546825; "N" = 4.640475; "O" = 0.09335776776421496; "P" = 0.09839297464651384; "Q" = 4.92312014495; "R" = 44.30808130455; "S" = 0.01306194909227805; "T" = 0.01403707709036024 }
    8 = @{ "I" = 0.1399128257358017; "J" = 0.1426634080409682; "O" = 0.3226929381796237; "P" = 0.3400972286002662; "S" = 0.04514888082569952; "T" = 0.04851942969740222 }
    9 = @{ "I" = 0.1399128257358017; "J" = 0.1426634080409682; "M" = 3.750299666666667; "N" = 11.250899; "O" = 0.2263472631531553; "P" = 0.238555195331833; "Q" = 11.93617625689133; "R" = 107.425586312022; "S" = 0.03166888518532306; "T" = 0.03403309717191817 }
    10 = @{ "I" = 0.1399128257358017; "J" = 0.1426634080409682; "M" = 2.543696; "N" = 5.087391999999999; "O" = 0.1535233658822184; "P" = 0.1078690504900635; "Q" = 8.095887395296; "R" = 48.575324371776; "S" = 0.02147988793705254; "T" = 0.01538896636505572 }
    11 = @{ "I" = 0.1399128257358017; "J" = 0.1426634080409682; "M" = 3.381336; "N" = 10.144008; "O" = 0.2040786650207874; "P" = 0.2150855509313235; "Q" = 10.761866001936; "R" = 96.85679401742399; "S" = 0.02855322269544847; "T" = 0.03068483771623185 }
    12 = @{ "G" = 6.716415; "H" = 20.149245; "I" = 0.295254005988679; "J" = 0.3010584805972866; "M" = 1.546825; "N" = 4.640475; "O" = 0.09335776776421496; "P" = 0.09839297464651384; "Q" = 10.389118632375; "R" = 93.502067691375; "S" = 0.02756425492254523; "T" = 0.02962203944852679 }
    13 = @{ "G" = 6.716415; "H" = 20.149245; "I" = 0.295254005988679; "J" = 0.3010584805972866; "O" = 0.3226929381796237; "P" = 0.3400972286002662; "Q" = 35.91019040906001; "R" = 323.19171368154; "S" = 0.09527638270179106; "T" = 0.1023891548977442 }
    14 = @{ "G" = 6.716415; "H" = 20.149245; "I" = 0.295254005988679; "J" = 0.3010584805972866; "M" = 3.750299666666667; "N" = 11.250899; "O" = 0.2263472631531553; "P" = 0.238555195331833; "Q" = 25.188568935695; "R" = 226.697120421255; "S" = 0.06682993619054281; "T" = 0.07181906464519056 }
    15 = @{ "G" = 6.716415; "H" = 20.149245; "I" = 0.295254005988679; "J" = 0.3010584805972866; "M" = 2.543696; "N" = 5.087391999999999; "O" = 0.1535233658822184; "P" = 0.1078690504900635; "Q" = 17.08451796984; "R" = 102.50710781904; "S" = 0.04532838878959067; "T" = 0.03247489244401049 }
    16 = @{ "G" = 6.716415; "H" = 20.149245; "I" = 0.295254005988679; "J" = 0.3010584805972866; "M" = 3.381336; "N" = 10.144008; "O" = 0.2040786650207874; "P" = 0.2150855509313235; "Q" = 22.71045583044; "R" = 204.39410247396; "S" = 0.06025504338420919; "T" = 0.06475332916181455 }
    17 = @{ "G" = 1.315755; "H" = 2.63151; "I" = 0.05784066866767975; "J" = 0.03931851552137887; "M" = 1.546825; "N" = 4.640475; "O" = 0.09335776776421496; "P" = 0.09839297464651384; "Q" = 2.035242727875; "R" = 12.21145636725; "S" = 0.005399875712804151; "T" = 0.003868665700833592 }
    18 = @{ "G" = 1.315755; "H" = 2.63151; "I" = 0.05784066866767975; "J" = 0.03931851552137887; "O" = 0.3226929381796237; "P" = 0.3400972286002662; "Q" = 7.03485603282; "R" = 42.20913619692001; "S" = 0.01866477531864768; "T" = 0.0133721181614975 }
    19 = @{ "G" = 1.315755; "H" = 2.63151; "I" = 0.05784066866767975; "J" = 0.03931851552137887; "M" = 3.750299666666667; "N" = 11.250899; "O" = 0.2263472631531553; "P" = 0.238555195331833; "Q" = 4.934475537915; "R" = 29.60685322749; "S" = 0.01309207705187777; "T" = 0.009379636150360245 }
    20 = @{ "G" = 1.315755; "H" = 2.63151; "I" = 0.05784066866767975; "J" = 0.03931851552137887; "M" = 2.543696; "N" = 5.087391999999999; "O" = 0.1535233658822184; "P" = 0.1078690504900635; "Q" = 3.34688073048; "R" = 13.38752292192; "S" = 0.008879894138740364; "T" = 0.004241250935969961 }
    21 = @{ "G" = 1.315755; "H" = 2.63151; "I" = 0.05784066866767975; "J" = 0.03931851552137887; "M" = 3.381336; "N" = 10.144008; "O" = 0.2040786650207874; "P" = 0.2150855509313235; "Q" = 4.44900974868; "R" = 26.69405849208; "S" = 0.01180404644560977; "T" = 0.008456844572717568 }
    22 = @{ "G" = 3.501345333333334; "H" = 10.504036; "I" = 0.1539193507274987; "J" = 0.1569452909178086; "M" = 1.546825; "N" = 4.640475; "O" = 0.09335776776421496; "P" = 0.09839297464651384; "Q" = 5.415968495233334; "R" = 48.74371645710001; "S" = 0.01436956699963658; "T" = 0.01544231403016568 }
    23 = @{ "G" = 3.501345333333334; "H" = 10.504036; "I" = 0.1539193507274987; "J" = 0.1569452909178086; "O" = 0.3226929381796237; "P" = 0.3400972286002662; "Q" = 18.72040033379022; "R" = 168.483603004112; "S" = 0.04966868752895657; "T" = 0.05337665848300922 }
    24 = @{ "G" = 3.501345333333334; "H" = 10.504036; "I" = 0.1539193507274987; "J" = 0.1569452909178086; "M" = 3.750299666666667; "N" = 11.250899; "O" = 0.2263472631531553; "P" = 0.238555195331833; "Q" = 13.13109423648489; "R" = 118.179848128364; "S" = 0.03483922378347996; "T" = 0.03744011453130918 }
    25 = @{ "G" = 3.501345333333334; "H" = 10.504036; "I" = 0.1539193507274987; "J" = 0.1569452909178086; "M" = 2.543696; "N" = 5.087391999999999; "O" = 0.1535233658822184; "P" = 0.1078690504900635; "Q" = 8.906358119018666; "R" = 53.438148714112; "S" = 0.02363021679809128; "T" = 0.01692953951019079 }
    26 = @{ "G" = 3.501345333333334; "H" = 10.504036; "I" = 0.1539193507274987; "J" = 0.1569452909178086; "M" = 3.381336; "N" = 10.144008; "O" = 0.2040786650207874; "P" = 0.2150855509313235; "Q" = 11.839225024032; "R" = 106.553025216288; "S" = 0.0314116556173343; "T" = 0.07743387511742585 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$row").Value = $cols[$col]
    }
}
